$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing value for week 12 (B13)
$ws.Range("B13").Value = 680

# Add new row for week 13
$ws.Range("A14").Value = 13
$ws.Range("B14").Value = 478
